# Generate Report for Handoff
# Adds two new "Ready for handoff" file entries (row3 = .md file, row4 = .png file)
# to the Overview sheet and to each locale sheet (zh-cn, de-de), and refreshes the
# existing row2 entry to point at a new source file.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# New / updated source file identifiers
# ---------------------------------------------------------------------------
$row2Name = "33d09cfd-4cd1-456e-b953-d60d016afa83.png"
$row3Name = "ada99b3e-4f92-49b1-a9f1-fb686f92f019.md"
$row4Name = "b535d7ad-5a8c-4bb0-816f-0c790b887877.png"

$row2Target_zhcn = "e2108d0630f5625cfdcd063fd471e6f0a73653c5.png"
$row2Target_dede = "e2108d0630f5625cfdcd063fd471e6f0a73653c5.png"
$row3Target_zhcn = "ada99b3e-4f92-49b1-a9f1-fb686f92f019.f222d8a22d1a23dcb9fa06093c897b802e203882.zh-cn.xlf"
$row3Target_dede = "ada99b3e-4f92-49b1-a9f1-fb686f92f019.f222d8a22d1a23dcb9fa06093c897b802e203882.de-de.xlf"
$row4Target_zhcn = "8ea207edca619ab1cb915bf6ace504de82f1acd6.png"
$row4Target_dede = "8ea207edca619ab1cb915bf6ace504de82f1acd6.png"

$statusText = "Ready for handoff"
$newHandoffDate = "2016-48-20 06:48:00"

$zhcnDateTime = "2016-03-20 06:47:57"
$dedeDateTime = "2016-03-20 06:48:00"
$epoch = "0001-01-01 00:00:00"

$row2Url = "https://github.com/OpenLocalizationTest/oltest/blob/218164d3f7d9a3c15347d1bfce4087ea74bdb2a8/e2e/$row2Name"
$row3Url = "https://github.com/OpenLocalizationTest/oltest/blob/218164d3f7d9a3c15347d1bfce4087ea74bdb2a8/e2e/$row3Name"
$row4Url = "https://github.com/OpenLocalizationTest/oltest/blob/218164d3f7d9a3c15347d1bfce4087ea74bdb2a8/e2e/$row4Name"

$row2UrlTarget_zhcn = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e932d778e72fa928677436a8e4e524686e6c6c7a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$row2Target_zhcn"
$row3UrlTarget_zhcn = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e932d778e72fa928677436a8e4e524686e6c6c7a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$row3Target_zhcn"
$row4UrlTarget_zhcn = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e932d778e72fa928677436a8e4e524686e6c6c7a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$row4Target_zhcn"

$row2UrlTarget_dede = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7f61ed617cbac77a17fae113d0ded1ebf778eeb6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$row2Target_dede"
$row3UrlTarget_dede = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7f61ed617cbac77a17fae113d0ded1ebf778eeb6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$row3Target_dede"
$row4UrlTarget_dede = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7f61ed617cbac77a17fae113d0ded1ebf778eeb6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$row4Target_dede"

# ---------------------------------------------------------------------------
# Overview sheet: update row 2, add rows 3 and 4
# ---------------------------------------------------------------------------
$overview.Range("A2").Value = $row2Name
$overview.Range("B2").Value = $statusText
$overview.Range("C2").Value = $statusText
$overview.Range("D2").Value = $newHandoffDate

$overview.Range("A3").Value = $row3Name
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText
$overview.Range("D3").Value = $newHandoffDate

$overview.Range("A4").Value = $row4Name
$overview.Range("B4").Value = $statusText
$overview.Range("C4").Value = $statusText
$overview.Range("D4").Value = $newHandoffDate

# Update the existing A2 hyperlink in-place (keep its rId / position).
foreach ($h in $overview.Hyperlinks) {
    $h.Address = $row2Url
    $h.TextToDisplay = $row2Name
}

# Add the two new hyperlinks for the freshly-added rows.
$overview.Hyperlinks.Add($overview.Range("A3"), $row3Url, "", "", $row3Name) | Out-Null
$overview.Hyperlinks.Add($overview.Range("A4"), $row4Url, "", "", $row4Name) | Out-Null

# ---------------------------------------------------------------------------
# Helper to populate a locale sheet (zh-cn / de-de)
# ---------------------------------------------------------------------------
function Set-LocaleSheet($ws, $dateTimeText, $row2Target, $row3Target, $row4Target, $row2UrlTarget, $row3UrlTarget, $row4UrlTarget) {
    # --- row 2 (existing row, refreshed) ---
    $ws.Range("A2").Value = $row2Name
    $ws.Range("B2").Value = ".png"
    $ws.Range("C2").Value = $statusText
    $ws.Range("D2").Value = $row2Target
    $ws.Range("E2").Value = $dateTimeText
    $ws.Range("H2").Value = $epoch
    $ws.Range("I2").Value = "IsDependency"
    $ws.Range("J2").Value = "e2e\" + $row3Name

    # --- row 3 (new row) ---
    $ws.Range("A3").Value = $row3Name
    $ws.Range("A3").Style = "Hyperlink"
    $ws.Range("B3").Value = ".md"
    $ws.Range("B3").Style = "Hyperlink"
    $ws.Range("C3").Value = $statusText
    $ws.Range("D3").Value = $row3Target
    $ws.Range("D3").Style = "Hyperlink"
    $ws.Range("E3").Value = $dateTimeText
    $ws.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("H3").Value = $epoch
    $ws.Range("I3").Value = "Include"

    # --- row 4 (new row) ---
    $ws.Range("A4").Value = $row4Name
    $ws.Range("A4").Style = "Hyperlink"
    $ws.Range("B4").Value = ".png"
    $ws.Range("B4").Style = "Hyperlink"
    $ws.Range("C4").Value = $statusText
    $ws.Range("D4").Value = $row4Target
    $ws.Range("D4").Style = "Hyperlink"
    $ws.Range("E4").Value = $dateTimeText
    $ws.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("H4").Value = $epoch
    $ws.Range("I4").Value = "IsDependency"
    $ws.Range("J4").Value = "e2e\" + $row3Name

    # Update the 3 existing hyperlinks (A2, B2, D2) in place.
    $existing = @()
    foreach ($h in $ws.Hyperlinks) {
        $existing += $h
    }
    $existing[0].Address = $row2Url
    $existing[0].TextToDisplay = $row2Name
    $existing[1].Address = $row2Url
    $existing[1].TextToDisplay = ".png"
    $existing[2].Address = $row2UrlTarget
    $existing[2].TextToDisplay = $row2Target

    # Add the new hyperlinks for rows 3 and 4.
    $ws.Hyperlinks.Add($ws.Range("A3"), $row3Url, "", "", $row3Name) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B3"), $row3Url, "", "", ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D3"), $row3UrlTarget, "", "", $row3Target) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A4"), $row4Url, "", "", $row4Name) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B4"), $row4Url, "", "", ".png") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D4"), $row4UrlTarget, "", "", $row4Target) | Out-Null
}

Set-LocaleSheet $zhcn $zhcnDateTime $row2Target_zhcn $row3Target_zhcn $row4Target_zhcn $row2UrlTarget_zhcn $row3UrlTarget_zhcn $row4UrlTarget_zhcn
Set-LocaleSheet $dede $dedeDateTime $row2Target_dede $row3Target_dede $row4Target_dede $row2UrlTarget_dede $row3UrlTarget_dede $row4UrlTarget_dede

Write-Output "Report generated for handoff"
